$d = $word.ActiveDocument

# Locate the paragraph that ends with "... Editora: ATLAS" (the last
# bibliography entry, "Gestão Empresarial - Estratégias Organizacionais...").
# Immediately after it come three paragraphs that must be removed:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# The empty paragraph that follows those three (right before the
# page-break paragraph) must be kept.

$marker = "Gestão Empresarial - Estratégias Organizacionais Autor: Bertero, C. O. Editora: ATLAS"

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $marker) {
        $found = $true

        # The next three paragraphs (empty / "Ver no Jupiter..." /
        # "© 2020 ...") are deleted by removing the Range that spans
        # from the start of the paragraph right after the marker to the
        # end of the third paragraph following it.
        $startPara = $d.Paragraphs.Item($i + 1)
        $endPara = $d.Paragraphs.Item($i + 3)

        $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $deleteRange.Delete()

        break
    }
}

if (-not $found) {
    throw "Could not locate bibliography marker paragraph"
}
